$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Fix Random Seed for benchmarking"
$ws.Range("A8").Value = "Remove Bye from matches in optimization"

$ws.Range("A8").Select()
